# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H/I/J/K/L/M/N columns across several
# Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 15874968
$ws.Range("I100").Value = 23810412
$ws.Range("J100").Value = 4080.7144
$ws.Range("K100").Value = 23810412
$ws.Range("L100").Value = 4080.7144
$ws.Range("M100").Value = -23809871
$ws.Range("N100").Value = -5162.7144

# Row 136
$ws.Range("H136").Value = 41900
$ws.Range("J136").Value = 41900
$ws.Range("L136").Value = 41900
$ws.Range("N136").Value = -52100

# Row 139
$ws.Range("H139").Value = 78050
$ws.Range("J139").Value = 78050
$ws.Range("L139").Value = 78050
$ws.Range("N139").Value = -88330

# Row 140
$ws.Range("H140").Value = 98000
$ws.Range("J140").Value = 98000
$ws.Range("L140").Value = 98000
$ws.Range("N140").Value = -108360

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 956919.0600000001
$ws.Range("I32").Value = 11190.881
$ws.Range("K32").Value = 11190.881
$ws.Range("M32").Value = -10903.881

# Row 61
$ws.Range("H61").Value = 25993.1
$ws.Range("I61").Value = 5236.25
$ws.Range("J61").Value = 109020.5
$ws.Range("K61").Value = 5236.25
$ws.Range("L61").Value = 109020.5
$ws.Range("M61").Value = -5024.25
$ws.Range("N61").Value = -109444.5

# Row 74
$ws.Range("H74").Value = 694.16327
$ws.Range("I74").Value = 673.3261
$ws.Range("J74").Value = 1013.6667
$ws.Range("K74").Value = 673.3261
$ws.Range("L74").Value = 1013.6667
$ws.Range("M74").Value = 200.6739
$ws.Range("N74").Value = -2761.6667

# Row 77
$ws.Range("H77").Value = 694.16327
$ws.Range("I77").Value = 673.3261
$ws.Range("J77").Value = 1013.6667
$ws.Range("K77").Value = 3366.6305
$ws.Range("L77").Value = 5068.3335
$ws.Range("M77").Value = 1001.3695
$ws.Range("N77").Value = -13804.3335

# Row 136
$ws.Range("H136").Value = 25993.1
$ws.Range("I136").Value = 5236.25
$ws.Range("J136").Value = 109020.5
$ws.Range("K136").Value = 15708.75
$ws.Range("L136").Value = 327061.5
$ws.Range("M136").Value = -13158.75
$ws.Range("N136").Value = -332161.5

# Row 138
$ws.Range("H138").Value = 69189.86
$ws.Range("J138").Value = 69189.86
$ws.Range("L138").Value = 69189.86
$ws.Range("N138").Value = -79469.86

# Row 139
$ws.Range("H139").Value = 38283.75
$ws.Range("J139").Value = 38283.75
$ws.Range("L139").Value = 38283.75
$ws.Range("N139").Value = -48563.75

# Row 140
$ws.Range("H140").Value = 108381
$ws.Range("J140").Value = 108381
$ws.Range("L140").Value = 108381
$ws.Range("N140").Value = -118741

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4796
$ws.Range("I86").Value = 1997.5
$ws.Range("J86").Value = 15990
$ws.Range("K86").Value = 1997.5
$ws.Range("L86").Value = 15990
$ws.Range("M86").Value = -874.5
$ws.Range("N86").Value = -18236

# Row 89
$ws.Range("H89").Value = 4796
$ws.Range("I89").Value = 1997.5
$ws.Range("J89").Value = 15990
$ws.Range("K89").Value = 9987.5
$ws.Range("L89").Value = 79950
$ws.Range("M89").Value = -4371.5
$ws.Range("N89").Value = -91182

# Row 105
$ws.Range("H105").Value = 2124.8
$ws.Range("I105").Value = 2163.5908
$ws.Range("K105").Value = 2163.5908
$ws.Range("M105").Value = -416.5907999999999

# Row 134
$ws.Range("H134").Value = 6001.4814
$ws.Range("I134").Value = 604.7143
$ws.Range("K134").Value = 1814.1429
$ws.Range("M134").Value = 720.8571000000002

# Row 138
$ws.Range("H138").Value = 50775
$ws.Range("J138").Value = 50775
$ws.Range("L138").Value = 50775
$ws.Range("N138").Value = -61055

# Row 140
$ws.Range("H140").Value = 86780
$ws.Range("J140").Value = 86780
$ws.Range("L140").Value = 86780
$ws.Range("N140").Value = -97140

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9697.671
$ws.Range("I31").Value = 3301.4211
$ws.Range("K31").Value = 3301.4211
$ws.Range("M31").Value = -3006.4211

# Row 34
$ws.Range("H34").Value = 9697.671
$ws.Range("I34").Value = 3301.4211
$ws.Range("K34").Value = 3301.4211
$ws.Range("M34").Value = -3099.4211

# Row 86
$ws.Range("H86").Value = 32934.2
$ws.Range("I86").Value = 6127.9375
$ws.Range("K86").Value = 6127.9375
$ws.Range("M86").Value = -5004.9375

# Row 89
$ws.Range("H89").Value = 32934.2
$ws.Range("I89").Value = 6127.9375
$ws.Range("K89").Value = 30639.6875
$ws.Range("M89").Value = -25023.6875

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1109.1923
$ws.Range("I5").Value = 1048.6666
$ws.Range("J5").Value = 1191.7273
$ws.Range("K5").Value = 3145.9998
$ws.Range("L5").Value = 3575.1819
$ws.Range("M5").Value = -3033.9998
$ws.Range("N5").Value = -3799.1819

# Row 38
$ws.Range("H38").Value = 305.28
$ws.Range("I38").Value = 648.25
$ws.Range("J38").Value = 143.88235
$ws.Range("K38").Value = 1944.75
$ws.Range("L38").Value = 431.64705
$ws.Range("M38").Value = -1597.75
$ws.Range("N38").Value = -1125.64705

# Row 86
$ws.Range("H86").Value = 187.5
$ws.Range("I86").Value = 187.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 562.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 623.5
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 187.5
$ws.Range("I89").Value = 187.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1687.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4240.5
$ws.Range("N89").ClearContents()

# Row 121
$ws.Range("H121").Value = 88975.95
$ws.Range("I121").Value = 3284.8
$ws.Range("J121").Value = 154892.23
$ws.Range("K121").Value = 9854.400000000001
$ws.Range("L121").Value = 464676.6900000001
$ws.Range("M121").Value = -8544.400000000001
$ws.Range("N121").Value = -467296.6900000001

# Row 135
$ws.Range("H135").Value = 1109.1923
$ws.Range("I135").Value = 1048.6666
$ws.Range("J135").Value = 1191.7273
$ws.Range("K135").Value = 9437.999400000001
$ws.Range("L135").Value = 10725.5457
$ws.Range("M135").Value = -6902.999400000001
$ws.Range("N135").Value = -15795.5457

# Row 137
$ws.Range("H137").Value = 4050.5454
$ws.Range("I137").Value = 2089.8572
$ws.Range("J137").Value = 4965.533
$ws.Range("K137").Value = 6269.571599999999
$ws.Range("L137").Value = 14896.599
$ws.Range("M137").Value = -1169.571599999999
$ws.Range("N137").Value = -25096.599

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2943.6667
$ws.Range("I132").Value = 2452.5334
$ws.Range("J132").Value = 5399.3335
$ws.Range("K132").Value = 7357.600199999999
$ws.Range("L132").Value = 16198.0005
$ws.Range("M132").Value = -4827.600199999999
$ws.Range("N132").Value = -21258.0005

# Row 133
$ws.Range("H133").Value = 54000
$ws.Range("J133").Value = 54000
$ws.Range("L133").Value = 54000
$ws.Range("N133").Value = -64120

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 140
$ws.Range("H140").Value = 75779.5
$ws.Range("J140").Value = 75779.5
$ws.Range("L140").Value = 75779.5
$ws.Range("N140").Value = -86139.5

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# Row 138
$ws.Range("H138").Value = 61350.777
$ws.Range("J138").Value = 61350.777
$ws.Range("L138").Value = 61350.777
$ws.Range("N138").Value = -71630.777

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1213.1034
$ws.Range("I136").Value = 765
$ws.Range("J136").Value = 2621.4285
$ws.Range("K136").Value = 2295
$ws.Range("L136").Value = 7864.2855
$ws.Range("M136").Value = 255
$ws.Range("N136").Value = -12964.2855

# Row 138
$ws.Range("H138").Value = 139800
$ws.Range("J138").Value = 139800
$ws.Range("L138").Value = 139800
$ws.Range("N138").Value = -150080
